# Colocacion de llenado automatico en gettypes de combo box con los datos que hay.
# Adds a new "Predecidos" column (Y) to the report sheet, mirroring the header
# style already used by the other header cells, and appends a new data row
# (row 5) with the latest game entry, including the computed "Predecidos"
# value in column Y.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell Y1 -----------------------------------------------
$ws.Range("Y1").Value = "Predecidos"
# Match the formatting already applied to the rest of the header row
# (bold font, thin border, centered/top aligned) by copying the format
# from the neighboring header cell.
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New data row (row 5) ----------------------------------------------
$ws.Range("A5").Value = "2024-07-28 16:00:48"
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("N5").Value = 10
$ws.Range("O5").Value = 6
$ws.Range("P5").Value = 3
$ws.Range("R5").Value = 5
$ws.Range("T5").Value = 20
$ws.Range("U5").Value = 0.6666666666666666
$ws.Range("V5").Value = "./Data/Electromecanica.xlsx"
$ws.Range("X5").Value = "No es Simulación"

# Newly computed "Predecidos" value for the row that was just played.
$ws.Range("Y5").Value = 9
